# ---------------------------------------------------------------------------
# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement worker/period table grows from 3 periods (2505-2507)
# to 4 periods (2505-2508) for the same 11 workers, so the detail grid grows
# from 33 data rows (16-48) to 44 data rows (16-59). The trailer block
# ("firma / nombre del representante legal") is pushed further down the
# sheet to make room, and the header totals (Valor Mora, Cant. Periodos)
# are refreshed to match.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($r, $tipoDoc, $numDoc, $nombre, $periodo, $valorMora, $salario)
    $ws.Cells.Item($r, 2).Value = $tipoDoc
    $ws.Cells.Item($r, 3).Value = $numDoc
    $ws.Cells.Item($r, 4).Value = $nombre
    $ws.Cells.Item($r, 5).Value = $periodo
    $ws.Cells.Item($r, 6).Value = $valorMora
    $ws.Cells.Item($r, 7).Value = $salario
}

# 1) Move the trailer block (old rows 53-54) out of the way first by
#    inserting 11 fresh rows right after the last existing data row (48).
#    That pushes the blank spacer rows + trailer rows down to 60-63 / 64-65,
#    exactly where they belong once the grid has 44 data rows (16-59).
$ws.Range("49:59").Insert() | Out-Null

# 2) The newly inserted rows don't carry the table's borders/number formats.
#    Clone them from the existing interior row (47) and, for the new final
#    row, from row 48 (heavier "bottom of table" border).
for ($r = 49; $r -le 58; $r++) {
    $ws.Range("B47:J47").Copy($ws.Range("B" + $r + ":J" + $r)) | Out-Null
}
$ws.Range("B48:J48").Copy($ws.Range("B59:J59")) | Out-Null

# 3) Re-write the full detail grid: 11 workers x 4 periods (2505..2508),
#    grouped by period.
Set-Row 16 "CC" "19615176" "ALVARO JOSE RAMIREZ MEDINA" "2505" 61880 1547000
Set-Row 17 "CC" "73156966" "DONALDO ENRIQUE RAMIREZ MEDINA" "2505" 61880 1547000
Set-Row 18 "CC" "45371797" "ADRIANA SOFIA RAMIREZ CASTRO" "2505" 61880 1547000
Set-Row 19 "CC" "92449681" "EDINALDO RUIZ CHIQUILLO" "2505" 56940 1423500
Set-Row 20 "CC" "45531043" "JOHANA ISABEL MONTALVO ACUNA" "2505" 61880 1547000
Set-Row 21 "CC" "9154526" "ULFRAN ENRIQUE PEREZ CASSIANI" "2505" 61880 1547000
Set-Row 22 "CC" "1002319738" "WILFRAN ENRIQUE PEREZ ZUÑIGA" "2505" 61880 1547000
Set-Row 23 "CC" "9042660" "JOSE GREGORIO CHAVES VERGARA" "2505" 56940 1423500
Set-Row 24 "CC" "1193333959" "LEIMER MEDRANO BERRIO" "2505" 56940 1423500
Set-Row 25 "CC" "9156656" "CARLOS ALEJANDRO TAPIA RAMIREZ" "2505" 61880 1547000
Set-Row 26 "CC" "1002322422" "EVER LUIS GALAN VALIENTE" "2505" 61880 1547000
Set-Row 27 "CC" "19615176" "ALVARO JOSE RAMIREZ MEDINA" "2506" 61880 1547000
Set-Row 28 "CC" "73156966" "DONALDO ENRIQUE RAMIREZ MEDINA" "2506" 61880 1547000
Set-Row 29 "CC" "45371797" "ADRIANA SOFIA RAMIREZ CASTRO" "2506" 61880 1547000
Set-Row 30 "CC" "92449681" "EDINALDO RUIZ CHIQUILLO" "2506" 56940 1423500
Set-Row 31 "CC" "45531043" "JOHANA ISABEL MONTALVO ACUNA" "2506" 61880 1547000
Set-Row 32 "CC" "9154526" "ULFRAN ENRIQUE PEREZ CASSIANI" "2506" 61880 1547000
Set-Row 33 "CC" "1002319738" "WILFRAN ENRIQUE PEREZ ZUÑIGA" "2506" 61880 1547000
Set-Row 34 "CC" "9042660" "JOSE GREGORIO CHAVES VERGARA" "2506" 56940 1423500
Set-Row 35 "CC" "1193333959" "LEIMER MEDRANO BERRIO" "2506" 56940 1423500
Set-Row 36 "CC" "9156656" "CARLOS ALEJANDRO TAPIA RAMIREZ" "2506" 61880 1547000
Set-Row 37 "CC" "1002322422" "EVER LUIS GALAN VALIENTE" "2506" 61880 1547000
Set-Row 38 "CC" "19615176" "ALVARO JOSE RAMIREZ MEDINA" "2507" 61880 1547000
Set-Row 39 "CC" "73156966" "DONALDO ENRIQUE RAMIREZ MEDINA" "2507" 61880 1547000
Set-Row 40 "CC" "45371797" "ADRIANA SOFIA RAMIREZ CASTRO" "2507" 61880 1547000
Set-Row 41 "CC" "92449681" "EDINALDO RUIZ CHIQUILLO" "2507" 56940 1423500
Set-Row 42 "CC" "45531043" "JOHANA ISABEL MONTALVO ACUNA" "2507" 61880 1547000
Set-Row 43 "CC" "9154526" "ULFRAN ENRIQUE PEREZ CASSIANI" "2507" 61880 1547000
Set-Row 44 "CC" "1002319738" "WILFRAN ENRIQUE PEREZ ZUÑIGA" "2507" 61880 1547000
Set-Row 45 "CC" "9042660" "JOSE GREGORIO CHAVES VERGARA" "2507" 56940 1423500
Set-Row 46 "CC" "1193333959" "LEIMER MEDRANO BERRIO" "2507" 56940 1423500
Set-Row 47 "CC" "9156656" "CARLOS ALEJANDRO TAPIA RAMIREZ" "2507" 61880 1547000
Set-Row 48 "CC" "1002322422" "EVER LUIS GALAN VALIENTE" "2507" 61880 1547000
Set-Row 49 "CC" "19615176" "ALVARO JOSE RAMIREZ MEDINA" "2508" 61880 1547000
Set-Row 50 "CC" "73156966" "DONALDO ENRIQUE RAMIREZ MEDINA" "2508" 61880 1547000
Set-Row 51 "CC" "45371797" "ADRIANA SOFIA RAMIREZ CASTRO" "2508" 61880 1547000
Set-Row 52 "CC" "92449681" "EDINALDO RUIZ CHIQUILLO" "2508" 56940 1423500
Set-Row 53 "CC" "45531043" "JOHANA ISABEL MONTALVO ACUNA" "2508" 61880 1547000
Set-Row 54 "CC" "9154526" "ULFRAN ENRIQUE PEREZ CASSIANI" "2508" 61880 1547000
Set-Row 55 "CC" "1002319738" "WILFRAN ENRIQUE PEREZ ZUÑIGA" "2508" 61880 1547000
Set-Row 56 "CC" "9042660" "JOSE GREGORIO CHAVES VERGARA" "2508" 56940 1423500
Set-Row 57 "CC" "1193333959" "LEIMER MEDRANO BERRIO" "2508" 56940 1423500
Set-Row 58 "CC" "9156656" "CARLOS ALEJANDRO TAPIA RAMIREZ" "2508" 61880 1547000
Set-Row 59 "CC" "1002322422" "EVER LUIS GALAN VALIENTE" "2508" 61880 1547000

# 4) Refresh the header totals: Valor Mora (sum over all periods) and
#    Cant. Periodos (now 4 instead of 3).
$ws.Range("E11").Value = 2663440
$ws.Range("F13").Value = 4

# 5) Re-merge the trailer block's cells at its new location (the insert
#    already unmerged them when it shifted the rows down).
$ws.Range("B64:C64").Merge() | Out-Null
$ws.Range("H64:J64").Merge() | Out-Null
$ws.Range("B65:C65").Merge() | Out-Null
$ws.Range("H65:J65").Merge() | Out-Null
